$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.152.61"
$ws.Range("E2").Value = "  -2.28%  "
$ws.Range("D3").Value = "1.647.53"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3926"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3870"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.64%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.372"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "49.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08578"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.109"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001284"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.525"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.71%  "
$ws.Range("D17").Value = "1.631.81"
$ws.Range("E17").Value = "  -7.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06919"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.940"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("E23").Value = "  -2.84%  "
$ws.Range("D24").Value = "24.143.76"
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.420"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.883"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "158.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.240"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.55%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "140.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.309"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.496"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("D33").Value = "1.824.75"
$ws.Range("E33").Value = "  -6.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08149"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.820"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02928"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9664"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2693"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.09164"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.425"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7539"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "13.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6914"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.462"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.098"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08386"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.270"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "133.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.89%  "
